# Append the new W/t (disorder-scan) data rows produced by the new
# simulation run / analysis script, mirroring rows 2-16 above them.
#
# Columns: A=File, B=nu, C=c_crit, D=expansion orders, E=Picture

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 17; A = "offdiagE0W10.txt";         B = "1.142 [.121, .160]"; C = 0.297; D = 3121; E = "E0W10.png" },
    @{ Row = 18; A = "offdiagE0W12.txt";         B = "1.159 [.130, .185]"; C = 0.432; D = 3121; E = "E0W12.png" },
    @{ Row = 19; A = "offdiagE0W14.txt";         B = "1.27 [.225, .297]";  C = 0.58;  D = 3121; E = "E0W14.png" },
    @{ Row = 20; A = "offdiagE0W16.txt";         B = "1.396 [.339, .464]"; C = 0.76;  D = 3121; E = "E0W16.png" },
    @{ Row = 21; A = "offdiagE0W8.txt";          B = "1.133 [.088, .172]"; C = 0.169; D = 3121; E = "E0W8.png" },
    @{ Row = 22; A = "offdiagE2W8.txt";          B = "1.129 [.106, .145]"; C = 0.186; D = 3121; E = "E2W8.png" },
    @{ Row = 23; A = "offdiagE4W8.txt";          B = "1.158 [.143, .176]"; C = 0.288; D = 3121; E = "E4W8.png" },
    @{ Row = 24; A = "offdiagE6W16take2.txt";    B = "1.892 [.769, 2.04]"; C = 0.884; D = 3121; E = "E6W16take2.png" },
    @{ Row = 25; A = "offdiagE6W14take2.txt";    B = "1.366 [.339, .402]"; C = 0.744; D = 3121; E = "E6W14take2.png" },
    @{ Row = 26; A = "offdiagE4W16take2.txt";    B = "1.408 [.346, .487]"; C = 0.789; D = 3121; E = "E4W16take2.png" },
    @{ Row = 27; A = "offdiagE4W14take2.txt";    B = "1.273 [.236, .311]"; C = 0.621; D = 3121; E = "E4W14take2.png" },
    @{ Row = 28; A = "offdiagE2W16take2.txt";    B = "1.406 [.344, .471]"; C = 0.766; D = 3121; E = "E2W16take2.png" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}

# Reflect the author's final cursor position/scroll recorded in the saved view.
$ws.Range("A29").Select()
